$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Approved" column (I), shifting it to J
$ws.Columns("I").Insert()

# Header for the new column
$ws.Range("I1").Value = "Distribution channel code"

# Data values for the new column
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"

# Match the column width used in the target workbook (closest representable value)
$ws.Columns("I").ColumnWidth = 21.7142857
